$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns keep text formatting so numeric-looking
# strings (e.g. "0.9997", "1.000") are not coerced into numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '30.094.07'
$ws.Range("E2").Value = '  -1.24%  '
$ws.Range("D3").Value = '1.851.35'
$ws.Range("E3").Value = '  -3.09%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").Value = '233.55'
$ws.Range("E5").Value = '  -2.24%  '
$ws.Range("D6").Value = '0.9997'
$ws.Range("E6").Value = '  +0.05%  '
$ws.Range("D7").Value = '0.4679'
$ws.Range("E7").Value = '  -1.86%  '
$ws.Range("D8").Value = '0.2818'
$ws.Range("E8").Value = '  -0.67%  '
$ws.Range("D9").Value = '0.06558'
$ws.Range("E9").Value = '  -2.11%  '
$ws.Range("D10").Value = '20.27'
$ws.Range("E10").Value = '  +4.22%  '
$ws.Range("D11").Value = '0.07778'
$ws.Range("E11").Value = '  +0.46%  '
$ws.Range("D12").Value = '97.16'
$ws.Range("E12").Value = '  -5.98%  '
$ws.Range("D13").Value = '1.847.06'
$ws.Range("E13").Value = '  -3.34%  '
$ws.Range("D14").Value = '5.079'
$ws.Range("E14").Value = '  -2.28%  '
$ws.Range("D15").Value = '0.6683'
$ws.Range("E15").Value = '  -0.07%  '
$ws.Range("D16").Value = '285.05'
$ws.Range("E16").Value = '  +2.70%  '
$ws.Range("D17").Value = '30.112.93'
$ws.Range("E17").Value = '  -1.22%  '
$ws.Range("D18").Value = '0.9999'
$ws.Range("E18").Value = '  +0.11%  '
$ws.Range("D19").Value = '12.58'
$ws.Range("E19").Value = '  -0.57%  '
$ws.Range("B20").Value = 'Uniswap'
$ws.Range("C20").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D20").Value = '5.352'
$ws.Range("E20").Value = '  -0.81%  '
$ws.Range("B21").Value = 'ShibaInu'
$ws.Range("C21").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D21").Value = '0.000007243'
$ws.Range("E21").Value = '  -3.17%  '
$ws.Range("B22").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C22").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D22").Value = '2.089.97'
$ws.Range("E22").Value = '  -3.07%  '
$ws.Range("B23").Value = 'BinanceUSD'
$ws.Range("C23").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D23").Value = '0.9999'
$ws.Range("E23").Value = '  +0.03%  '
$ws.Range("D24").Value = '6.131'
$ws.Range("E24").Value = '  -2.45%  '
$ws.Range("D25").Value = '167.72'
$ws.Range("E25").Value = '  +0.84%  '
$ws.Range("D26").Value = '9.289'
$ws.Range("E26").Value = '  -0.56%  '
$ws.Range("D27").Value = '19.09'
$ws.Range("E27").Value = '  -0.85%  '
$ws.Range("D28").Value = '1.932'
$ws.Range("E28").Value = '  -6.64%  '
$ws.Range("D29").Value = '1.337'
$ws.Range("E29").Value = '  -3.16%  '
$ws.Range("D30").Value = '0.09740'
$ws.Range("E30").Value = '  -2.55%  '
$ws.Range("D31").Value = '4.395'
$ws.Range("E31").Value = '  -4.98%  '
$ws.Range("D32").Value = '1.467'
$ws.Range("E32").Value = '  -2.83%  '
$ws.Range("D33").Value = '4.084'
$ws.Range("E33").Value = '  -3.73%  '
$ws.Range("D34").Value = '0.04653'
$ws.Range("E34").Value = '  -0.63%  '
$ws.Range("D35").Value = '0.6992'
$ws.Range("E35").Value = '  -3.96%  '
$ws.Range("D36").Value = '1.082'
$ws.Range("E36").Value = '  -2.87%  '
$ws.Range("D37").Value = '0.9990'
$ws.Range("E37").Value = '  +0.08%  '
$ws.Range("D38").Value = '2.707'
$ws.Range("E38").Value = '  +0.09%  '
$ws.Range("D39").Value = '0.01855'
$ws.Range("E39").Value = '  -2.60%  '
$ws.Range("D40").Value = '6.330'
$ws.Range("E40").Value = '  +0.10%  '
$ws.Range("D41").Value = '2.500'
$ws.Range("E41").Value = '  -4.17%  '
$ws.Range("D42").Value = '71.59'
$ws.Range("E42").Value = '  -3.77%  '
$ws.Range("D43").Value = '0.8605'
$ws.Range("E43").Value = '  -0.04%  '
$ws.Range("D44").Value = '1.931'
$ws.Range("E44").Value = '  -1.41%  '
$ws.Range("B45").Value = 'Quant'
$ws.Range("C45").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D45").Value = '103.87'
$ws.Range("E45").Value = '  -1.80%  '
$ws.Range("B46").Value = 'PaxDollar'
$ws.Range("C46").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D46").Value = '0.9994'
$ws.Range("E46").Value = '  +0.04%  '
$ws.Range("D47").Value = '0.4159'
$ws.Range("E47").Value = '  -2.35%  '
$ws.Range("D48").Value = '1.010.08'
$ws.Range("E48").Value = '  +6.14%  '
$ws.Range("D49").Value = '7.251'
$ws.Range("E49").Value = '  -2.50%  '
$ws.Range("D50").Value = '9.209'
$ws.Range("E50").Value = '  +4.70%  '
$ws.Range("D51").Value = '33.72'
$ws.Range("E51").Value = '  -2.67%  '
